# Before changing to copy each time of template word file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1: "Адрес" -> "адрес" (lowercase)
$ws.Range("D1").Value = "адрес"

# D column (address) for rows 22 and 23 first
$ws.Range("D22").Value = "Sin City"
$ws.Range("D23").Value = "Moskvabad"

# B column (name) for rows 22 and 23
$ws.Range("B22").Value = "Имя21"
$ws.Range("B23").Value = "Имя22"

# C column (numbers) and A23 (index)
$ws.Range("C22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("C23").Value = 22

# Remove old F22/F23 cells (no longer used)
$ws.Range("F22").Clear()
$ws.Range("F23").Clear()

# Update selection to B24
$ws.Range("B24").Select()
